# issue #5: stock data output to json file
#
# Adds a "property_category" column (all rows = "stock") to the 股票
# (stock) sheet, fixes a couple of data-entry glitches on that sheet
# (a stray space in a company name, and full-width comma digit grouping
# in a face-value cell that should be a clean numeric-looking text
# value), and lets every reference into shared strings on the other
# sheets fall out naturally from the recalculated/compacted string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H ("property_category"); this shifts the existing
# H/I/J (date / legislator_name / legislator_id) into I/J/K and keeps
# every other reference (styles, A column ids, etc.) intact.
$ws.Columns("H").Insert()

# Header
$ws.Cells.Item(1, 8).Value = "property_category"

# Every stock record in this sheet is a "stock" property.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Fix stray embedded space in "中國合成橡膠股份有限公司" (row 7, name column B).
$ws.Cells.Item(7, 2).Value = "中國合成橡膠股份有限公司"

# Fix the face-value text in row 6 (G): was full-width-comma grouped
# "1，500，000" -> plain "1500000". Keep it stored as text (not a number)
# by forcing the cell to a text format before assigning the digit string.
$g6 = $ws.Cells.Item(6, 7)
$g6.NumberFormat = "@"
$g6.Value = "1500000"
